# Update protocol docs sheet (Sheet1):
#  - 'sack'  -> 'sack, seqNum, sackedList'
#  - 'reset' -> 'reset, unicodeReason'
# (this also causes the now-unused 'reset'/'sack' shared strings to be
#  dropped and the two new strings appended at the end of the shared
#  string table, matching the upstream OOXML diff)
# Also move the sheet selection from A24:F24 to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters for shared-string append order: 'sack' string must be
# added before the 'reset' string so they land at indices 45 and 46
# respectively (matching the target sharedStrings.xml ordering).
$ws.Range("A12").Value = "sack, seqNum, sackedList"
$ws.Range("A11").Value = "reset, unicodeReason"

$ws.Range("A11").Select()
